$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything from E1 down to I4 (old trailing columns / rows no longer used)
$ws.Range("A1:I4").ClearContents()

# New header row
$ws.Range("A1").Value = "Qtd_Nós"
$ws.Range("B1").Value = "Ativos"
$ws.Range("C1").Value = "Distancia"
$ws.Range("D1").Value = "Tempo"

# New data row
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = 5344
$ws.Range("D2").Value = 1319.392914772034
